$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $text) {
    $rng = $ws.Range($cellAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "60.528.86"
$ws.Range("E2").Value = "  +1.52%  "
$ws.Range("D3").Value = "2.606.77"
$ws.Range("E3").Value = "  +1.96%  "
$ws.Range("E4").Value = "  -0.22%  "
Set-TextValue "D5" "514.95"
$ws.Range("E5").Value = "  +2.75%  "
Set-TextValue "D6" "153.80"
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("E7").Value = "  +0.07%  "
Set-TextValue "D8" "0.600"
$ws.Range("E8").Value = "  +5.14%  "
$ws.Range("D9").Value = "2.618.00"
$ws.Range("E9").Value = "  +1.33%  "
Set-TextValue "D10" "6.70"
$ws.Range("E10").Value = "  +3.81%  "
$ws.Range("E11").Value = "  +1.22%  "
$ws.Range("E12").Value = "  +2.29%  "
Set-TextValue "D13" "0.129"
$ws.Range("E13").Value = "  +2.09%  "
$ws.Range("D14").Value = "3.063.54"
$ws.Range("E14").Value = "  +1.36%  "
$ws.Range("D15").Value = "60.537.35"
$ws.Range("E15").Value = "  +1.21%  "
Set-TextValue "D16" "21.67"
$ws.Range("E16").Value = "  +0.93%  "
$ws.Range("E17").Value = "  +2.11%  "
$ws.Range("D18").Value = "2.610.25"
$ws.Range("E18").Value = "  +0.78%  "
Set-TextValue "D19" "4.76"
$ws.Range("E19").Value = "  +0.36%  "
Set-TextValue "D20" "357.97"
$ws.Range("E20").Value = "  +5.95%  "
Set-TextValue "D21" "10.63"
$ws.Range("E21").Value = "  +3.37%  "
Set-TextValue "D22" "6.21"
$ws.Range("E22").Value = "  +2.94%  "
Set-TextValue "D23" "1.00"
$ws.Range("E23").Value = "  +0.22%  "
Set-TextValue "D24" "61.06"
$ws.Range("E24").Value = "  +2.51%  "
Set-TextValue "D25" "0.426"
$ws.Range("E25").Value = "  +2.01%  "
$ws.Range("D26").Value = "2.725.16"
$ws.Range("E26").Value = "  +1.06%  "
$ws.Range("E27").Value = "  +1.51%  "
Set-TextValue "D28" "0.997"
$ws.Range("E28").Value = "  +0.21%  "
$ws.Range("E29").Value = "  +0.00%  "
Set-TextValue "D30" "7.34"
$ws.Range("E30").Value = "  -0.67%  "
Set-TextValue "D31" "1.00"
$ws.Range("E31").Value = "  +0.02%  "
Set-TextValue "D32" "19.46"
$ws.Range("E32").Value = "  +1.82%  "
Set-TextValue "D33" "1.59"
$ws.Range("E33").Value = "  +2.62%  "
Set-TextValue "D34" "5.91"
$ws.Range("E34").Value = "  +4.29%  "
Set-TextValue "D35" "150.49"
$ws.Range("E35").Value = "  -3.36%  "
Set-TextValue "D36" "4.02"
$ws.Range("E36").Value = "  +2.15%  "
$ws.Range("E37").Value = "  +0.34%  "
Set-TextValue "D38" "0.907"
$ws.Range("E38").Value = "  +7.02%  "
$ws.Range("E39").Value = "  +1.60%  "
Set-TextValue "D40" "0.845"
$ws.Range("E40").Value = "  +1.28%  "
Set-TextValue "D41" "36.25"
$ws.Range("E41").Value = "  +2.43%  "
Set-TextValue "D42" "3.75"
$ws.Range("E42").Value = "  +0.35%  "
Set-TextValue "D43" "290.14"
$ws.Range("E43").Value = "  -2.32%  "
$ws.Range("E44").Value = "  +2.84%  "
$ws.Range("E45").Value = "  +1.04%  "
Set-TextValue "D48" "19.65"
$ws.Range("E48").Value = "  +0.67%  "
Set-TextValue "D49" "4.97"
$ws.Range("E49").Value = "  +1.87%  "
$ws.Range("E50").Value = "  +1.96%  "
$ws.Range("E51").Value = "  +0.44%  "

# Row 46/47 swap: Hedera <-> FirstDigitalUSD
$ws.Range("B46").Value = "FirstDigitalUSD"
$ws.Range("C46").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue "D46" "0.997"
$ws.Range("E46").Value = "  +0.35%  "

$ws.Range("B47").Value = "Hedera"
$ws.Range("C47").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D47" "0.0557"
$ws.Range("E47").Value = "  -1.21%  "
